# Update cryptocurrency price/volume data per the latest symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value is numeric-looking
# text (Price/Volume columns). Numeric-looking text needs a leading apostrophe so
# Excel keeps storing it as literal text instead of converting it to a number.
$changes = @(
    @{ Cell = "D2"; Value = '246.49'; NumericText = $true }
    @{ Cell = "D3"; Value = '29.61'; NumericText = $true }
    @{ Cell = "E3"; Value = '9.81%'; NumericText = $true }
    @{ Cell = "D4"; Value = '5.165'; NumericText = $true }
    @{ Cell = "E4"; Value = '1.87%'; NumericText = $true }
    @{ Cell = "D5"; Value = '0.05708'; NumericText = $true }
    @{ Cell = "D6"; Value = '6.601'; NumericText = $true }
    @{ Cell = "D7"; Value = '3.073'; NumericText = $true }
    @{ Cell = "E7"; Value = '2.25%'; NumericText = $true }
    @{ Cell = "D8"; Value = '0.8558'; NumericText = $true }
    @{ Cell = "E8"; Value = '4.39%'; NumericText = $true }
    @{ Cell = "D9"; Value = '0.8692'; NumericText = $true }
    @{ Cell = "E9"; Value = '3.21%'; NumericText = $true }
    @{ Cell = "B10"; Value = 'WazirX'; NumericText = $false }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; NumericText = $false }
    @{ Cell = "D10"; Value = '0.1367'; NumericText = $true }
    @{ Cell = "E10"; Value = '2.88%'; NumericText = $true }
    @{ Cell = "B11"; Value = 'MandalaExchangeToken'; NumericText = $false }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; NumericText = $false }
    @{ Cell = "D11"; Value = '0.07076'; NumericText = $true }
    @{ Cell = "E11"; Value = '2.71%'; NumericText = $true }
    @{ Cell = "B12"; Value = 'BitrueCoin'; NumericText = $false }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; NumericText = $false }
    @{ Cell = "D12"; Value = '0.02924'; NumericText = $true }
    @{ Cell = "E12"; Value = '2.21%'; NumericText = $true }
    @{ Cell = "B13"; Value = 'BitMartToken'; NumericText = $false }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; NumericText = $false }
    @{ Cell = "D13"; Value = '0.09390'; NumericText = $true }
    @{ Cell = "E13"; Value = '-0.03%'; NumericText = $true }
    @{ Cell = "B14"; Value = 'BitForexToken'; NumericText = $false }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; NumericText = $false }
    @{ Cell = "D14"; Value = '0.001518'; NumericText = $true }
    @{ Cell = "E14"; Value = '-0.23%'; NumericText = $true }
    @{ Cell = "B15"; Value = 'CoinExToken'; NumericText = $false }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; NumericText = $false }
    @{ Cell = "D15"; Value = '0.04180'; NumericText = $true }
    @{ Cell = "E15"; Value = '1.96%'; NumericText = $true }
    @{ Cell = "B16"; Value = 'One'; NumericText = $false }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; NumericText = $false }
    @{ Cell = "D16"; Value = '0.0006019'; NumericText = $true }
    @{ Cell = "E16"; Value = '0.55%'; NumericText = $true }
    @{ Cell = "D17"; Value = '0.006121'; NumericText = $true }
    @{ Cell = "E17"; Value = '-0.35%'; NumericText = $true }
    @{ Cell = "D19"; Value = '3.485'; NumericText = $true }
    @{ Cell = "E19"; Value = '-0.70%'; NumericText = $true }
    @{ Cell = "D20"; Value = '2.187'; NumericText = $true }
    @{ Cell = "D22"; Value = '0.03323'; NumericText = $true }
    @{ Cell = "E22"; Value = '4.00%'; NumericText = $true }
    @{ Cell = "E23"; Value = '0.33%'; NumericText = $true }
    @{ Cell = "D24"; Value = '3.470'; NumericText = $true }
    @{ Cell = "E24"; Value = '-2.37%'; NumericText = $true }
    @{ Cell = "E25"; Value = '0.49%'; NumericText = $true }
    @{ Cell = "D26"; Value = '0.005028'; NumericText = $true }
    @{ Cell = "E26"; Value = '26.78%'; NumericText = $true }
    @{ Cell = "E27"; Value = '0.54%'; NumericText = $true }
    @{ Cell = "E28"; Value = '23.56%'; NumericText = $true }
    @{ Cell = "D40"; Value = '0.03743'; NumericText = $true }
    @{ Cell = "E40"; Value = '1.31%'; NumericText = $true }
    @{ Cell = "B41"; Value = 'BKEXToken'; NumericText = $false }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; NumericText = $false }
    @{ Cell = "D41"; Value = '0.1073'; NumericText = $true }
    @{ Cell = "E41"; Value = '1.44%'; NumericText = $true }
    @{ Cell = "B42"; Value = 'CEJI'; NumericText = $false }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; NumericText = $false }
    @{ Cell = "D42"; Value = '0.002540'; NumericText = $true }
    @{ Cell = "E42"; Value = '2.91%'; NumericText = $true }
    @{ Cell = "B43"; Value = 'KickToken'; NumericText = $false }
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; NumericText = $false }
    @{ Cell = "D43"; Value = '0.003476'; NumericText = $true }
    @{ Cell = "E43"; Value = '-36.76%'; NumericText = $true }
    @{ Cell = "D44"; Value = '0.009824'; NumericText = $true }
    @{ Cell = "E44"; Value = '4.46%'; NumericText = $true }
    @{ Cell = "D45"; Value = '0.00005213'; NumericText = $true }
    @{ Cell = "E45"; Value = '0.09%'; NumericText = $true }
    @{ Cell = "E46"; Value = '0.02%'; NumericText = $true }
    @{ Cell = "D47"; Value = '0.05799'; NumericText = $true }
    @{ Cell = "E47"; Value = '-42.85%'; NumericText = $true }
    @{ Cell = "D48"; Value = '0.002574'; NumericText = $true }
    @{ Cell = "E48"; Value = '-0.84%'; NumericText = $true }
    @{ Cell = "E49"; Value = '0.02%'; NumericText = $true }
    @{ Cell = "E50"; Value = '0.02%'; NumericText = $true }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    if ($change.NumericText) {
        # Leading apostrophe forces Excel to keep the value as literal text
        $range.Value = "'" + $change.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $change.Value
    }
}

Write-Output "Updated $($changes.Count) cells"
